# Regenerate the s_vals data to filter save games.
# Updates columns B, C, D, E, G for rows 2-8 with the newly computed values.
# Column F (Win flag) is left unchanged. Column G is the sum of B:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 4.23731228292506)
    3 = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 10.19245300693656, 35.69412875252057)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    5 = @(0.003208871385164791, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.075923540317991)
    6 = @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.406728370586922)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    8 = @(3.286832544864788, 250555.8564151394, 0.7527432677738641, 10.19245300693656, 250570.0884439589)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
